# Regenerate save_data column G ("K") using new computation (K instead of Strike#),
# recalc std/mean and write the new s_vals (K column) back to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4,4,4,1,1,3,0,3,4,4,6,5,4,5,6,2,5,6,6,8,4,2,4,3,6,6,4,6,2,4,4,4,2,1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
